$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($Row, $F, $H, $I, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $T, $U, $V)
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 8).Value = $H
    $ws.Cells.Item($Row, 9).Value = $I
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $N
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = $Q
    $ws.Cells.Item($Row, 18).Value = $R
    $ws.Cells.Item($Row, 19).Value = $S
    $ws.Cells.Item($Row, 20).Value = $T
    $ws.Cells.Item($Row, 21).Value = $U
    $ws.Cells.Item($Row, 22).Value = $V
}

Set-RowData 2 82 "living_rooms" "distractor" "f" "stimuli/img_xu1p3.png" 75.27659574468085 56.68085106382978 65.97872340425532 47 7 7 7 6 6 6
Set-RowData 3 83 "bedrooms" "distractor" "f" "stimuli/img_71mhq.png" 69.34210526315789 47.02631578947368 58.18421052631579 38 5 5 5 5 5 5
Set-RowData 4 84 "kitchens" "target" "j" "stimuli/img_as3da.png" 84.53125 63 73.765625 32 9 9 9 9 9 9
Set-RowData 5 85 "kitchens" "target" "j" "stimuli/img_0j24m.png" 63.6969696969697 35.75757575757576 49.72727272727273 33 3 3 3 3 3 3
Set-RowData 6 86 "living_rooms" "distractor" "f" "stimuli/img_abobq.png" 75.18421052631579 54.13157894736842 64.65789473684211 38 6 6 6 6 6 6
Set-RowData 7 87 "bedrooms" "distractor" "f" "stimuli/img_088sa.png" 67.36363636363636 42.12121212121212 54.74242424242424 33 4 4 4 4 4 4
Set-RowData 8 88 "bedrooms" "distractor" "f" "stimuli/img_huisn.png" 73.63888888888889 46.36111111111111 60 36 5 5 5 5 5 5
Set-RowData 9 89 "kitchens" "target" "j" "stimuli/img_uspja.png" 54.90909090909091 29.12121212121212 42.01515151515152 33 2 2 2 2 2 2
Set-RowData 10 90 "bedrooms" "distractor" "f" "stimuli/img_h13c3.png" 71.80555555555556 47.86111111111111 59.83333333333334 36 5 5 5 5 5 5
Set-RowData 11 91 "bedrooms" "distractor" "f" "stimuli/img_d3t0o.png" 66.95121951219512 42.92682926829269 54.9390243902439 41 4 4 4 4 4 4
Set-RowData 12 92 "kitchens" "target" "j" "stimuli/img_qmand.png" 86.11764705882354 71.02941176470588 78.57352941176471 34 10 10 10 10 10 10
Set-RowData 13 93 "living_rooms" "distractor" "f" "stimuli/img_02alv.png" 61.8 37.8 49.8 45 4 4 4 4 5 4
Set-RowData 14 94 "bedrooms" "distractor" "f" "stimuli/img_b89t4.png" 71.41463414634147 47.85365853658536 59.63414634146342 41 5 5 5 5 5 5
Set-RowData 15 95 "kitchens" "target" "j" "stimuli/img_u9f9l.png" 77.78571428571429 57.25 67.51785714285714 28 7 7 7 7 7 7
Set-RowData 16 96 "living_rooms" "distractor" "f" "stimuli/img_qrc78.png" 76.2 59.875 68.03749999999999 40 7 7 7 7 6 7
Set-RowData 17 97 "kitchens" "target" "j" "stimuli/img_oz18d.png" 78.93939393939394 61.03030303030303 69.98484848484848 33 8 8 8 8 8 8
Set-RowData 18 98 "bedrooms" "distractor" "f" "stimuli/img_bdz92.png" 63.72222222222222 42.63888888888889 53.18055555555556 36 4 4 4 4 4 4
Set-RowData 19 99 "bedrooms" "distractor" "f" "stimuli/img_mqnl6.png" 70.7560975609756 45.68292682926829 58.21951219512195 41 5 5 5 5 5 5
Set-RowData 20 100 "kitchens" "target" "j" "stimuli/img_z5osu.png" 71.42857142857143 47.34285714285714 59.38571428571429 35 5 5 5 5 5 5
Set-RowData 21 101 "living_rooms" "distractor" "f" "stimuli/img_g9od8.png" 59.34883720930232 37.83720930232558 48.59302325581395 43 4 4 4 4 4 4
Set-RowData 22 102 "kitchens" "target" "j" "stimuli/img_t1cr9.png" 73.66666666666667 53.51515151515152 63.59090909090909 33 6 6 6 6 6 6
Set-RowData 23 103 "kitchens" "target" "j" "stimuli/img_8fpog.png" 85.41666666666667 72.30555555555556 78.86111111111111 36 10 10 10 10 10 10
Set-RowData 24 104 "kitchens" "target" "j" "stimuli/img_h1yyu.png" 64.85294117647059 46.61764705882353 55.73529411764706 34 4 4 4 4 4 4
Set-RowData 25 105 "kitchens" "target" "j" "stimuli/img_l9t30.png" 67.2 43.14285714285715 55.17142857142858 35 4 4 4 4 4 4
Set-RowData 26 106 "kitchens" "target" "j" "stimuli/img_qbkdt.png" 69.45714285714286 50.02857142857143 59.74285714285715 35 5 5 5 5 5 5
Set-RowData 27 107 "kitchens" "target" "j" "stimuli/img_hfz8w.png" 55.46153846153846 27.28205128205128 41.37179487179487 39 2 2 2 2 2 2
Set-RowData 28 108 "living_rooms" "distractor" "f" "stimuli/img_4974k.png" 79.89130434782609 58.19565217391305 69.04347826086956 46 7 7 7 7 7 7
Set-RowData 29 109 "kitchens" "target" "j" "stimuli/img_k3abb.png" 35.54054054054054 16.54054054054054 26.04054054054054 37 1 1 1 1 1 1
Set-RowData 30 110 "living_rooms" "distractor" "f" "stimuli/img_fmgjx.png" 79.90000000000001 56.975 68.4375 40 7 7 7 7 7 7
Set-RowData 31 111 "living_rooms" "distractor" "f" "stimuli/img_zh8ms.png" 59.82608695652174 39.43478260869565 49.6304347826087 46 4 4 4 4 4 4
Set-RowData 32 112 "bedrooms" "distractor" "f" "stimuli/img_uttnz.png" 69.91891891891892 49.91891891891892 59.91891891891892 37 5 5 5 5 5 5
Set-RowData 33 113 "kitchens" "target" "j" "stimuli/img_57os5.png" 82.70588235294117 65.73529411764706 74.22058823529412 34 9 9 9 9 9 9
Set-RowData 34 114 "kitchens" "target" "j" "stimuli/img_01w8b.png" 78.91891891891892 61.21621621621622 70.06756756756756 37 8 8 8 8 8 8
Set-RowData 35 115 "kitchens" "target" "j" "stimuli/img_5949k.png" 60.8 39.2 50 35 3 3 3 3 3 3
Set-RowData 36 116 "kitchens" "target" "j" "stimuli/img_ensho.png" 72.7948717948718 54.56410256410256 63.67948717948718 39 6 6 6 6 6 6
Set-RowData 37 117 "living_rooms" "distractor" "f" "stimuli/img_lpas9.png" 59.36585365853659 39.09756097560975 49.23170731707317 41 4 4 4 4 4 4
Set-RowData 38 118 "bedrooms" "distractor" "f" "stimuli/img_twj5p.png" 67.71739130434783 42.08695652173913 54.90217391304348 46 4 4 4 4 5 4
Set-RowData 39 119 "kitchens" "target" "j" "stimuli/img_fhm45.png" 76.75 57.71875 67.234375 32 7 7 7 7 7 7
Set-RowData 40 120 "kitchens" "target" "j" "stimuli/img_faly8.png" 33.41176470588236 19.23529411764706 26.32352941176471 34 1 1 1 1 1 1
Set-RowData 41 121 "living_rooms" "distractor" "f" "stimuli/img_lpr0l.png" 77.04651162790698 59.86046511627907 68.45348837209303 43 7 7 7 7 7 7
